$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the amount in T2 (was 103842, now 93098)
$ws.Range("T2").Value = 93098

# Leave the active selection on T2 (was T3)
$ws.Range("T2").Select()
